# Add two new diary entries (2/13 and 2/16) below the existing entries,
# in the previously-blank rows 21 and 22 of the diary worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: reuse the formatting already used by the entry directly above
# it (row 20) - same column styles: date / highlighted N/A participants /
# goal / achievements / reflection / mood.
$ws.Range("A20:G20").Copy()
$ws.Range("A21:G21").PasteSpecial(-4122)

# --- Row 22: reuse the formatting used by an earlier plain entry (row 15)
# for columns B:G, then fix column A's number format to the short date
# style ("d-mmm") used throughout the log by copying it from row 19.
$ws.Range("A15:G15").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)
$ws.Range("A19").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 21 values (Thursday 2/13 class entry) ---
$ws.Range("A21").Value = 43874
$ws.Range("B21").Value = "5:00 -7:50 pm"
$ws.Range("C21").Value = "N/A"
$ws.Range("D21").Value = "big picture and more key expert practices"
$ws.Range("E21").Value = "We learned about various stakeholders important to the development of software and how experts work along different levels of abstrction and how they prioritize work"
$ws.Range("F21").Value = "I really found the key expert practice ""do something else"" very helpful. Often I get stuck on trying to understand how a piece of code works and I will fixate on it. I think I needed to hear that it's good practice to stop and search other areas of the code or to do somethin entirely different while your mind sorts out hte problem. "
$ws.Range("G21").Value = "Good"

# --- Row 22 values (Sunday 2/16 group-work entry) ---
$ws.Range("A22").Value = 43877
$ws.Range("B22").Value = "3:00 - 7:30pm"
$ws.Range("C22").Value = "Chris, Jay, Rafael"
$ws.Range("D22").Value = "Worked on finding the stakeholders for latest project"
$ws.Range("E22").Value = "we were able to accomplish our goal of finding stakeholders by searching thorugh forums, github, and documentation"
$ws.Range("F22").Value = "I learned the value of documentation, forums, and github. We were lucky to have such an engaged community with our OS project. By searching through forums we learned that our OS project is actively engaged with its users by implementnig features and reporting/fixing bugs. It was great to see such an active role from the developr side and makes me appreciate how much care goes into someones software"
# G22 is left blank, matching the source edit.

# --- Row heights to fit the newly-wrapped text ---
$ws.Rows("21").RowHeight = 153
$ws.Rows("22").RowHeight = 204

# --- Update the saved view/selection to match where the author ended up:
# scrolled down a couple rows, with the new F22 cell selected.
$ws.Range("F22").Select()
$excel.ActiveWindow.ScrollRow = 19
